# intersection_setup_file.xlsx: add Detector_Length / Detector_2_Length /
# Detect_Length / Detect_2_Length columns to support e2 (area) detectors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns in the middle of the table -----------------
# Insert before old column I (Sumo_Detector_2) -> new I holds Detector_Length,
# everything from old I onward shifts one column right (I..Q -> J..R).
$ws.Columns("I:I").Insert()
# Insert before the (now shifted) old column L (Detector_2_Distance, now at
# M) -> new M holds Detector_2_Length, everything from there shifts right
# again (M..R -> N..S).
$ws.Columns("M:M").Insert()

# --- Header row --------------------------------------------------------
# I1/M1 inherit the bold/centered/bordered header style automatically from
# the Insert() above (Excel copies the left-neighbour column's formatting).
$ws.Range("I1").Value = "Detector_Length"
$ws.Range("M1").Value = "Detector_2_Length"

# T1/U1 are brand-new trailing columns, so they start out with the default
# style; give them the same header formatting as the rest of row 1 by
# copying an existing header cell's format onto them.
$ws.Range("T1").Value = "Detect_Length"
$ws.Range("U1").Value = "Detect_2_Length"
$ws.Range("H1").Copy()
$ws.Range("T1:U1").PasteSpecial(-4122)

# --- Data rows -----------------------------------------------------------
# For each row: Detector_Distance (H) changes for rows that now represent
# "area" (e2) detectors, Detector_Length (I) gets the length of the new
# detector, Detect_Actual_Edge/Lane (N/P) keep the values shifted in from the
# old L/N columns (already correct after the column inserts above) except
# row 11 whose matched edge changed, Detect_Actual_Distance (R) gets the new
# measured distance, and Detect_Length (T) mirrors the Detector_Length.

$rows = @(
    @{ Row = 2;  H = 100; I = 10; N = "Bryce_NB_0";  P = "Bryce_NB";   R = -30.5;               T = 10 },
    @{ Row = 3;  H = 100; I = 10; N = "Bryce_NB_1";  P = "Bryce_NB";   R = -30.5;               T = 10 },
    @{ Row = 4;  H = 30;  I = 30; N = "Bryce_NB_2";  P = "Bryce_NB";   R = -9.1;                T = 30 },
    @{ Row = 5;  H = 100; I = 10; N = "Campus_EB_0"; P = "Campus_EB";  R = -30.5;               T = 10 },
    @{ Row = 6;  H = 30;  I = 30; N = "Campus_EB_1"; P = "Campus_EB";  R = -9.1;                T = 30 },
    @{ Row = 7;  H = 100; I = 10; N = "Bryce_SB_0";  P = "Bryce_SB";   R = -30.5;               T = 10 },
    @{ Row = 8;  H = 100; I = 10; N = "Bryce_SB_1";  P = "Bryce_SB";   R = -30.5;               T = 10 },
    @{ Row = 9;  H = 30;  I = 30; N = "Bryce_SB_2";  P = "Bryce_SB";   R = -9.1;                T = 30 },
    @{ Row = 11; H = 100; I = 10; N = "gneE1.93_1";  P = "gneE1.93";   R = -9.259999024640031;  T = 10 },
    @{ Row = 12; H = 30;  I = 30; N = "Campus_WB_2"; P = "Campus_WB";  R = -9.1;                T = 30 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("N$n").Value = $r.N
    $ws.Range("P$n").Value = $r.P
    $ws.Range("R$n").Value = $r.R
    $ws.Range("T$n").Value = $r.T
}

# Row 10 has no Sumo_Detector/RW_Detector/Detector_Distance data at all (it
# stays entirely blank there), and its N/O/P/Q values are already correct
# after the column shift above, so nothing else to do for it.
